# Atualizando dados com jogos de 07/04/2019
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 98-102 previously empty placeholder rows (s="3" on A/B only).
# Copy the formatting (incl. the bordered "s=10" style) used by the C
# column on the row directly above, so the newly written C98:C102
# cells pick up the same look instead of the engine's default style.
$ws.Range("C97").Copy() | Out-Null
$ws.Range("C98:C102").PasteSpecial(-4122) | Out-Null

# Game results added on 07/04/2019 (ISG beat AK three times, then AK
# beat ISG twice), each entry worth -1 in the C column.
$ws.Range("A98").Value = "ISG"
$ws.Range("B98").Value = "AK"
$ws.Range("C98").Value = -1

$ws.Range("A99").Value = "ISG"
$ws.Range("B99").Value = "AK"
$ws.Range("C99").Value = -1

$ws.Range("A100").Value = "ISG"
$ws.Range("B100").Value = "AK"
$ws.Range("C100").Value = -1

$ws.Range("A101").Value = "AK"
$ws.Range("B101").Value = "ISG"
$ws.Range("C101").Value = -1

$ws.Range("A102").Value = "AK"
$ws.Range("B102").Value = "ISG"
$ws.Range("C102").Value = -1

# Author left the selection on F102 (was F95) after entering the data.
$ws.Range("F102").Select() | Out-Null
